$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B6 becomes a true numeric value instead of a text value
$ws.Range("B6").Value = 1000135120

# Add new row 7 with the negotiation log entry (2025-10-15 22:31:31 update)
$ws.Range("A7").Value = "2025-10-15 17:31:30"

# Cedula must stay textual (leading-zero-safe) even though it looks numeric
$ws.Range("B7").Value = "'1000127336"
$ws.Range("B7").ClearFormats()

$ws.Range("C7").Value = "Paula"
$ws.Range("D7").Value = "TARJETA DE CRÉDITO"
$ws.Range("E7").Value = "****4376"
$ws.Range("F7").Value = "REDIFERIDO SIN PAGO"
$ws.Range("G7").Value = "36 cuotas"
$ws.Range("H7").Value = "34.19.100.134"
$ws.Range("I7").Value = "The Dalles"
$ws.Range("J7").Value = "Oregon"
$ws.Range("K7").Value = "United States"
$ws.Range("L7").Value = "2025-10-15 17:31:30"
$ws.Range("M7").Value = "****4376"
$ws.Range("N7").Value = "34.19.100.134"
$ws.Range("O7").Value = ""
$ws.Range("P7").Value = ""
